# Scheduled runner update: refresh computed profit/cost columns (H:N) on
# several sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per latest market
# data snapshot. Applies exact cell values per row; cells not listed are
# intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 295.66666
$ws.Range("I4").Value = 131
$ws.Range("J4").Value = 748.5
$ws.Range("K4").Value = 131
$ws.Range("L4").Value = 748.5
$ws.Range("M4").Value = -17
$ws.Range("N4").Value = -976.5
$ws.Range("H125").Value = 1887.7858
$ws.Range("I125").Value = 1391.25
$ws.Range("J125").Value = 2549.8333
$ws.Range("K125").Value = 12521.25
$ws.Range("L125").Value = 22948.4997
$ws.Range("M125").Value = -10061.25
$ws.Range("N125").Value = -27868.4997
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 1200.9286
$ws.Range("I127").Value = 745.625
$ws.Range("J127").Value = 1808
$ws.Range("K127").Value = 2236.875
$ws.Range("L127").Value = 5424
$ws.Range("M127").Value = 2723.125
$ws.Range("N127").Value = -15344
$ws.Range("H128").Value = 48000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 48000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 48000
$ws.Range("N128").Value = -57960
$ws.Range("H129").Value = 143667.11
$ws.Range("I129").Value = 566
$ws.Range("J129").Value = 157082.84
$ws.Range("K129").Value = 1698
$ws.Range("L129").Value = 471248.52
$ws.Range("M129").Value = 3302
$ws.Range("N129").Value = -481248.52
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 3698.8
$ws.Range("I131").Value = 3547
$ws.Range("J131").Value = 3800
$ws.Range("K131").Value = 10641
$ws.Range("L131").Value = 11400
$ws.Range("M131").Value = -5601
$ws.Range("N131").Value = -21480
$ws.Range("H132").Value = 3082.0715
$ws.Range("I132").Value = 3334.9092
$ws.Range("J132").Value = 2155
$ws.Range("K132").Value = 10004.7276
$ws.Range("L132").Value = 6465
$ws.Range("M132").Value = -7474.7276
$ws.Range("N132").Value = -11525
$ws.Range("H133").Value = 52780
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 52780
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 52780
$ws.Range("N133").Value = -62900
$ws.Range("H134").Value = 50000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 50000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H135").Value = 7870.4375
$ws.Range("I135").Value = 870
$ws.Range("J135").Value = 16871
$ws.Range("K135").Value = 7830
$ws.Range("L135").Value = 151839
$ws.Range("M135").Value = -5295
$ws.Range("N135").Value = -156909
$ws.Range("H136").Value = 42280
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 42280
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 42280
$ws.Range("N136").Value = -52480
$ws.Range("H137").Value = 68703.2
$ws.Range("I137").Value = 2068.8
$ws.Range("J137").Value = 102020.4
$ws.Range("K137").Value = 6206.400000000001
$ws.Range("L137").Value = 306061.2
$ws.Range("M137").Value = -3656.400000000001
$ws.Range("N137").Value = -311161.2
$ws.Range("H138").Value = 1549.75
$ws.Range("I138").Value = 555.2353000000001
$ws.Range("J138").Value = 3086.7273
$ws.Range("K138").Value = 1665.7059
$ws.Range("L138").Value = 9260.1819
$ws.Range("M138").Value = 3474.2941
$ws.Range("N138").Value = -19540.1819
$ws.Range("H139").Value = 52780
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 52780
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 52780
$ws.Range("N139").Value = -63060
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 2690
$ws.Range("I141").Value = 1619
$ws.Range("J141").Value = 3761
$ws.Range("K141").Value = 4857
$ws.Range("L141").Value = 11283
$ws.Range("M141").Value = 323
$ws.Range("N141").Value = -21643
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20877.25
$ws.Range("I32").Value = 20945.51
$ws.Range("K32").Value = 20945.51
$ws.Range("M32").Value = -20658.51
$ws.Range("H61").Value = 2616.3333
$ws.Range("I61").Value = 1664.9231
$ws.Range("J61").Value = 5090
$ws.Range("K61").Value = 1664.9231
$ws.Range("L61").Value = 5090
$ws.Range("M61").Value = -1452.9231
$ws.Range("N61").Value = -5514
$ws.Range("H74").Value = 2159.2222
$ws.Range("I74").Value = 2119.4
$ws.Range("K74").Value = 2119.4
$ws.Range("M74").Value = -1245.4
$ws.Range("H77").Value = 2159.2222
$ws.Range("I77").Value = 2119.4
$ws.Range("K77").Value = 10597
$ws.Range("M77").Value = -6229
$ws.Range("H97").Value = 777.5
$ws.Range("I97").Value = 847.5
$ws.Range("K97").Value = 847.5
$ws.Range("M97").Value = -351.5
$ws.Range("H136").Value = 2616.3333
$ws.Range("I136").Value = 1664.9231
$ws.Range("J136").Value = 5090
$ws.Range("K136").Value = 4994.7693
$ws.Range("L136").Value = 15270
$ws.Range("M136").Value = -2444.7693
$ws.Range("N136").Value = -20370
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 31392.086
$ws.Range("I134").Value = 38761.535
$ws.Range("J134").Value = 1914.2858
$ws.Range("K134").Value = 116284.605
$ws.Range("L134").Value = 5742.857400000001
$ws.Range("M134").Value = -113749.605
$ws.Range("N134").Value = -10812.8574
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -100
$ws.Range("H31").Value = 2826.4092
$ws.Range("I31").Value = 1434.5714
$ws.Range("K31").Value = 1434.5714
$ws.Range("M31").Value = -1139.5714
$ws.Range("H34").Value = 2826.4092
$ws.Range("I34").Value = 1434.5714
$ws.Range("K34").Value = 1434.5714
$ws.Range("M34").Value = -1232.5714
$ws.Range("H58").Value = 21936.084
$ws.Range("I58").Value = 1157.5454
$ws.Range("K58").Value = 1157.5454
$ws.Range("M58").Value = -954.5454
$ws.Range("H132").Value = 2502.4092
$ws.Range("I132").Value = 1243.0769
$ws.Range("K132").Value = 3729.2307
$ws.Range("M132").Value = -1199.2307
$ws.Range("H134").Value = 1066.875
$ws.Range("I134").Value = 883.3684
$ws.Range("K134").Value = 2650.1052
$ws.Range("M134").Value = -115.1052
$ws.Range("H136").Value = 21936.084
$ws.Range("I136").Value = 1157.5454
$ws.Range("K136").Value = 3472.6362
$ws.Range("M136").Value = -922.6361999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 240
$ws.Range("I13").Value = 250
$ws.Range("J13").Value = 233.33333
$ws.Range("K13").Value = 750
$ws.Range("L13").Value = 699.99999
$ws.Range("M13").Value = -582
$ws.Range("N13").Value = -1035.99999
$ws.Range("H22").Value = 9809.091
$ws.Range("J22").Value = 1150
$ws.Range("L22").Value = 3450
$ws.Range("N22").Value = -3788
$ws.Range("H27").Value = 9809.091
$ws.Range("J27").Value = 1150
$ws.Range("L27").Value = 3450
$ws.Range("N27").Value = -3654
$ws.Range("H129").Value = 2283
$ws.Range("I129").Value = 739.6
$ws.Range("J129").Value = 10000
$ws.Range("K129").Value = 2218.8
$ws.Range("L129").Value = 30000
$ws.Range("M129").Value = 2781.2
$ws.Range("N129").Value = -40000
$ws.Range("H131").Value = 763.34
$ws.Range("J131").Value = 780.6989
$ws.Range("L131").Value = 2342.0967
$ws.Range("N131").Value = -12422.0967
$ws.Range("H139").Value = 2075.0625
$ws.Range("I139").Value = 1262.6
$ws.Range("K139").Value = 3787.8
$ws.Range("M139").Value = 1352.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1399.3334
$ws.Range("I97").Value = 886.8
$ws.Range("J97").Value = 2424.4
$ws.Range("K97").Value = 886.8
$ws.Range("L97").Value = 2424.4
$ws.Range("M97").Value = -390.8
$ws.Range("N97").Value = -3416.4
$ws.Range("H118").Value = 41000
$ws.Range("J118").Value = 41000
$ws.Range("L118").Value = 41000
$ws.Range("N118").Value = -44314
$ws.Range("H126").Value = 5862.207
$ws.Range("I126").Value = 5125.2
$ws.Range("K126").Value = 15375.6
$ws.Range("M126").Value = -12905.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 311.82352
$ws.Range("J16").Value = 307
$ws.Range("L16").Value = 307
$ws.Range("N16").Value = -647
$ws.Range("H132").Value = 3131.5
$ws.Range("I132").Value = 3004
$ws.Range("K132").Value = 9012
$ws.Range("M132").Value = -6482
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3466.6667
$ws.Range("I96").Value = 1200
$ws.Range("J96").Value = 8000
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 8000
$ws.Range("M96").Value = 173
$ws.Range("N96").Value = -10746
$ws.Range("H132").Value = 1391.4
$ws.Range("J132").Value = 3124.5
$ws.Range("L132").Value = 9373.5
$ws.Range("N132").Value = -14433.5
$ws.Range("H136").Value = 21278122
$ws.Range("I136").Value = 34484148
$ws.Range("K136").Value = 103452444
$ws.Range("M136").Value = -103449894
